$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update timestamp message (row 1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 12:22"

# Update city names where ranking shifted (rows 13-16, 26-27)
$ws.Range("A13").Value = "Zaragoza"
$ws.Range("A14").Value = "Alacant/Alicante"
$ws.Range("A15").Value = "Araba/Alava"
$ws.Range("A16").Value = "Toledo"
$ws.Range("A26").Value = "Cantabria"
$ws.Range("A27").Value = "Granada"

# Update numeric stats (Casos totales, Casos activos, Recuperados, Muertes)
$ws.Range("B4").Value = 46587
$ws.Range("C4").Value = 24683
$ws.Range("D4").Value = 15626
$ws.Range("E4").Value = 6278
$ws.Range("B10").Value = 3969
$ws.Range("C10").Value = 650
$ws.Range("D10").Value = 3087
$ws.Range("E10").Value = 232
$ws.Range("B12").Value = 3279
$ws.Range("C12").Value = 1420
$ws.Range("D12").Value = 1640
$ws.Range("E12").Value = 219
$ws.Range("B13").Value = 3057
$ws.Range("C13").Value = 706
$ws.Range("D13").Value = 2017
$ws.Range("E13").Value = 334
$ws.Range("B14").Value = 2993
$ws.Range("C14").Value = 764
$ws.Range("D14").Value = 1902
$ws.Range("E14").Value = 327
$ws.Range("B15").Value = 2990
$ws.Range("C15").Value = 5092
$ws.Range("D15").Value = 4658
$ws.Range("E15").Value = 254
$ws.Range("B16").Value = 2984
$ws.Range("C16").Value = 2205
$ws.Range("D16").Value = 9768
$ws.Range("E16").Value = 403
$ws.Range("B26").Value = 1752
$ws.Range("C26").Value = 297
$ws.Range("D26").Value = 1345
$ws.Range("E26").Value = 110
$ws.Range("B27").Value = 1725
$ws.Range("C27").Value = 272
$ws.Range("D27").Value = 1296
$ws.Range("E27").Value = 157
$ws.Range("B47").Value = 502
$ws.Range("C47").Value = 103
$ws.Range("D47").Value = 332
$ws.Range("E47").Value = 67
$ws.Range("B49").Value = 464
$ws.Range("C49").Value = 115
$ws.Range("D49").Value = 301
$ws.Range("E49").Value = 48
